$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201, shifting existing rows 201:309 down to 202:310
$ws.Rows("201:201").Insert()

# Populate the newly inserted row 201 with the new data record
$ws.Cells.Item(201, 1).Value = 5
$ws.Cells.Item(201, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(201, 3).Value = "Maule"
$ws.Cells.Item(201, 4).Value = 44813
$ws.Cells.Item(201, 5).Value = 7
$ws.Cells.Item(201, 6).Value = 100112009
$ws.Cells.Item(201, 7).Value = "Acelga"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 1500
$ws.Cells.Item(201, 11).Value = 2500
$ws.Cells.Item(201, 12).Value = 2500
$ws.Cells.Item(201, 13).Value = 2500
$ws.Cells.Item(201, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(201, 15).Value = "Región del Maule"
$ws.Cells.Item(201, 16).Value = 625
$ws.Cells.Item(201, 17).Value = 4
$ws.Cells.Item(201, 18).Value = "Hortaliza"
